$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-7 with the new TPM-derived values (shifted sending-cluster
# mapping: the old "ECs" sending-cluster rows are dropped, and the
# remaining FAPs/MuSCs rows get refreshed numeric values).
$ws.Range("A2").Value2 = "FAPs"
$ws.Range("B2").Value2 = "Efnb3"
$ws.Range("C2").Value2 = "Epha4"
$ws.Range("D2").Value2 = "ECs"
$ws.Range("E2").Value2 = 2
$ws.Range("F2").Value2 = 0.6666666666666666
$ws.Range("G2").Value2 = 0.1498043333333333
$ws.Range("H2").Value2 = 0.449413
$ws.Range("I2").Value2 = 0.08722868471333377
$ws.Range("J2").Value2 = 0.08722868471333377
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 6.030956000000001
$ws.Range("N2").Value2 = 18.092868
$ws.Range("O2").Value2 = 0.364814105361131
$ws.Range("P2").Value2 = 0.3648141053611309
$ws.Range("Q2").Value2 = 0.9034633429426668
$ws.Range("R2").Value2 = 8.131170086484001
$ws.Range("S2").Value2 = 0.03182225457552302
$ws.Range("T2").Value2 = 0.03182225457552301
$ws.Range("A3").Value2 = "FAPs"
$ws.Range("B3").Value2 = "Efnb3"
$ws.Range("C3").Value2 = "Epha4"
$ws.Range("D3").Value2 = "FAPs"
$ws.Range("E3").Value2 = 2
$ws.Range("F3").Value2 = 0.6666666666666666
$ws.Range("G3").Value2 = 0.1498043333333333
$ws.Range("H3").Value2 = 0.449413
$ws.Range("I3").Value2 = 0.08722868471333377
$ws.Range("J3").Value2 = 0.08722868471333377
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 6.789877333333333
$ws.Range("N3").Value2 = 20.369632
$ws.Range("O3").Value2 = 0.4107214552505144
$ws.Range("P3").Value2 = 0.4107214552505143
$ws.Range("Q3").Value2 = 1.017153047335111
$ws.Range("R3").Value2 = 9.154377426016
$ws.Range("S3").Value2 = 0.03582669232504875
$ws.Range("T3").Value2 = 0.03582669232504874
$ws.Range("A4").Value2 = "FAPs"
$ws.Range("B4").Value2 = "Efnb3"
$ws.Range("C4").Value2 = "Epha4"
$ws.Range("D4").Value2 = "MuSCs"
$ws.Range("E4").Value2 = 2
$ws.Range("F4").Value2 = 0.6666666666666666
$ws.Range("G4").Value2 = 0.1498043333333333
$ws.Range("H4").Value2 = 0.449413
$ws.Range("I4").Value2 = 0.08722868471333377
$ws.Range("J4").Value2 = 0.08722868471333377
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 3.710753333333333
$ws.Range("N4").Value2 = 11.13226
$ws.Range("O4").Value2 = 0.2244644393883547
$ws.Range("P4").Value2 = 0.2244644393883547
$ws.Range("Q4").Value2 = 0.5558869292644444
$ws.Range("R4").Value2 = 5.002982363379999
$ws.Range("S4").Value2 = 0.01957973781276201
$ws.Range("T4").Value2 = 0.01957973781276201
$ws.Range("A5").Value2 = "MuSCs"
$ws.Range("B5").Value2 = "Efnb3"
$ws.Range("C5").Value2 = "Epha4"
$ws.Range("D5").Value2 = "ECs"
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 1.567570333333333
$ws.Range("H5").Value2 = 4.702711
$ws.Range("I5").Value2 = 0.9127713152866662
$ws.Range("J5").Value2 = 0.9127713152866662
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 6.030956000000001
$ws.Range("N5").Value2 = 18.092868
$ws.Range("O5").Value2 = 0.364814105361131
$ws.Range("P5").Value2 = 0.3648141053611309
$ws.Range("Q5").Value2 = 9.453947707238667
$ws.Range("R5").Value2 = 85.085529365148
$ws.Range("S5").Value2 = 0.3329918507856079
$ws.Range("T5").Value2 = 0.3329918507856079
$ws.Range("A6").Value2 = "MuSCs"
$ws.Range("B6").Value2 = "Efnb3"
$ws.Range("C6").Value2 = "Epha4"
$ws.Range("D6").Value2 = "FAPs"
$ws.Range("E6").Value2 = 3
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 1.567570333333333
$ws.Range("H6").Value2 = 4.702711
$ws.Range("I6").Value2 = 0.9127713152866662
$ws.Range("J6").Value2 = 0.9127713152866662
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 6.789877333333333
$ws.Range("N6").Value2 = 20.369632
$ws.Range("O6").Value2 = 0.4107214552505144
$ws.Range("P6").Value2 = 0.4107214552505143
$ws.Range("Q6").Value2 = 10.64361027470578
$ws.Range("R6").Value2 = 95.792492472352
$ws.Range("S6").Value2 = 0.3748947629254656
$ws.Range("T6").Value2 = 0.3748947629254655
$ws.Range("A7").Value2 = "MuSCs"
$ws.Range("B7").Value2 = "Efnb3"
$ws.Range("C7").Value2 = "Epha4"
$ws.Range("D7").Value2 = "MuSCs"
$ws.Range("E7").Value2 = 3
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 1.567570333333333
$ws.Range("H7").Value2 = 4.702711
$ws.Range("I7").Value2 = 0.9127713152866662
$ws.Range("J7").Value2 = 0.9127713152866662
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 3.710753333333333
$ws.Range("N7").Value2 = 11.13226
$ws.Range("O7").Value2 = 0.2244644393883547
$ws.Range("P7").Value2 = 0.2244644393883547
$ws.Range("Q7").Value2 = 5.816866839651111
$ws.Range("R7").Value2 = 52.35180155685999
$ws.Range("S7").Value2 = 0.2048847015755927
$ws.Range("T7").Value2 = 0.2048847015755927

# Remove the now-obsolete trailing rows (old rows 8-10, MuSCs x 3 target rows)
$ws.Range("A8:T10").ClearContents()
